$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Replace the author name "Manus AI" -> "Charles Kinyua Gitonga"
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute("Manus AI", $true, $false, $false, $false,
                                  $false, $true, 1, $false,
                                  "Charles Kinyua Gitonga", 2)

# ---------------------------------------------------------------------------
# 2. Locate the "Author:" paragraph (now ending in "Charles Kinyua Gitonga")
#    and the following "Date:" paragraph, then splice in the new
#    Student ID / Department / Institution / Email / Course / Instructor
#    fields (each "<b>Label:</b> value" pair separated by a manual line
#    break) between the author name and the date field, finally folding the
#    old "Date:" paragraph's runs into the same paragraph so everything
#    lives in a single paragraph as in the target document.
# ---------------------------------------------------------------------------
$authorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "Author:*") {
        $authorPara = $para
        break
    }
}

$datePara = $authorPara.Next()

# Collapsed insertion point just before the author paragraph's end-of-paragraph mark.
$insertAt = $authorPara.Range.End - 1
$ins = $d.Range($insertAt, $insertAt)

function Add-Field($label, $value) {
    $ins.InsertAfter([char]11)

    $labelStart = $ins.End
    $ins.InsertAfter($label)
    $labelEnd = $ins.End
    $d.Range($labelStart, $labelEnd).Bold = 1

    $ins.InsertAfter(" ")
    $ins.InsertAfter($value)
}

Add-Field "Student ID:" "SD23/77993/25"
Add-Field "Department:" "Computer Science"
Add-Field "Institution:" "Chuka University"
Add-Field "Email:" "cgkinyua@chuka.ac.ke"
Add-Field "Course:" "COSC 944 - Multi-Agent Systems"
Add-Field "Instructor:" "Prof. Marcel Odhiambo Ohanga"

# Line break before the Date field, then the Date field itself with the
# updated date value.
Add-Field "Date:" "November 29, 2025"

# Remove the now-duplicated old "Date:" paragraph entirely (including its
# paragraph mark), which leaves the author paragraph's own mark/style intact.
$datePara.Range.Delete()
